# Record a new lending: Member 1015 (Eesa Jesus) borrows Book 10005 (Ghafasi dar Magas)

$wb = $excel.ActiveWorkbook

$wsMembers  = $wb.Worksheets.Item("Members")
$wsBooks    = $wb.Worksheets.Item("Books")
$wsLendings = $wb.Worksheets.Item("Lendings")

# Members sheet: row 14 is member 1015 (Eesa Jesus) - record which book(s) they now have checked out
$wsMembers.Range("F14").Value = ", 10005"

# Books sheet: row 4 is book 10005 (Ghafasi dar Magas) - record which member borrowed it
$wsBooks.Range("E4").Value = 1015

# Lendings sheet: append a new lending record (book id, member id, date of lending)
$wsLendings.Range("C5").NumberFormat = "@"
$wsLendings.Range("A5").Value = 10005
$wsLendings.Range("B5").Value = 1015
$wsLendings.Range("C5").Value = "2024/7/17"
$wsLendings.Range("C5").NumberFormat = "General"
$wsLendings.Range("C5").Style = "Normal"

Write-Output "Lending recorded: Book 10005 -> Member 1015"
